# APResolution now retrieved from MovieDatabase.xlsx
# Insert a new "APResolution" column into Sheet1 right after the "StemLoop"
# column (i.e. before the old "Channel1" column), shifting every column
# from old-F onward one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at F; everything from F onward shifts right.
$ws.Columns.Item(6).Insert()

# Populate and style the new header cell to match its neighbours (bold,
# same font as the rest of row 1).
$ws.Cells.Item(1, 6).Value = "APResolution"
$ws.Cells.Item(1, 6).Font.Bold = $true

# Give the new column a sensible width similar to the other text columns
# (StemLoop / Comments use ~25.57).
$ws.Columns.Item(6).ColumnWidth = 25.5703125

# Rebuild the AutoFilter over the new, wider header range (drop it first so
# re-applying doesn't just toggle the existing filter off).
$ws.AutoFilterMode = $false
$ws.Range("A1:S2").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name tracks the full table dimension
# (one column beyond the AutoFilter's own range), so update it explicitly.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$T`$2"

# Restore the user's selection to cell F5 on the sheet.
$ws.Range("F5").Select() | Out-Null
